$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.847.05"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "3.134.66"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'597.69"
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").Value = "'139.11"
$ws.Range("E6").Value = "  -4.19%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.126.24"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "'34.38"
$ws.Range("E14").Value = "  -3.48%  "
$ws.Range("D15").Value = "3.649.21"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").Value = "63.808.36"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "3.133.62"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "'483.03"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").Value = "'87.88"
$ws.Range("E24").Value = "  +4.79%  "
$ws.Range("E25").Value = "  -5.79%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("E28").Value = "  -5.91%  "
$ws.Range("D29").Value = "'6.92"
$ws.Range("E30").Value = "  -3.07%  "
$ws.Range("D31").Value = "'27.10"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -6.85%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'52.54"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "0.0₃0736"
$ws.Range("E38").Value = "  -5.77%  "
$ws.Range("D39").Value = "'0.0392"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'427.43"
$ws.Range("E40").Value = "  -7.62%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.84"
$ws.Range("E41").Value = "  -10.70%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "2.872.91"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  -3.18%  "
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("E47").Value = "  -6.81%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D50").Value = "'25.51"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").Value = "'120.50"
$ws.Range("E51").Value = "  +0.59%  "
